# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 46074 to 46075 (i.e. bump the date by one day).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cur = $ws.Cells.Item($r, 3).Value2
    if ($cur -eq 46074) {
        $ws.Cells.Item($r, 3).Value = 46075
    }
}
